$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 17 data in source order)
$ws.Range("D2").Value = 44214
$ws.Range("M2").Value = 48

# Row 3 (was row 8 data in source order)
$ws.Range("D3").Value = 44585
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 6500
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6750
$ws.Range("S3").Value = 3375

# Row 4 (was row 10 data in source order)
$ws.Range("D4").Value = 45001
$ws.Range("M4").Value = 66
$ws.Range("N4").Value = 7500
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7773
$ws.Range("R4").Value = 'Provincia de Curicó'
$ws.Range("S4").Value = 3886

# Row 6 (was row 2 data in source order)
$ws.Range("D6").Value = 44614
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 6000
$ws.Range("O6").Value = 6000
$ws.Range("P6").Value = 6000
$ws.Range("R6").Value = 'Provincia de Linares'
$ws.Range("S6").Value = 3000

# Row 7 (was row 4 data in source order)
$ws.Range("D7").Value = 44627
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("S7").Value = 3000

# Row 8 (was row 19 data in source order)
$ws.Range("D8").Value = 44586
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 7000
$ws.Range("P8").Value = 7000
$ws.Range("S8").Value = 3500

# Row 9 (was row 12 data in source order)
$ws.Range("D9").Value = 44628
$ws.Range("M9").Value = 40
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6000
$ws.Range("P9").Value = 6000
$ws.Range("S9").Value = 3000

# Row 10 (was row 16 data in source order)
$ws.Range("D10").Value = 44588
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 6500
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 6750
$ws.Range("S10").Value = 3375

# Row 12 (was row 9 data in source order)
$ws.Range("D12").Value = 44606
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 7000
$ws.Range("O12").Value = 7000
$ws.Range("P12").Value = 7000
$ws.Range("S12").Value = 3500

# Row 13 (was row 15 data in source order)
$ws.Range("D13").Value = 44587
$ws.Range("M13").Value = 165
$ws.Range("N13").Value = 6500
$ws.Range("O13").Value = 7000
$ws.Range("P13").Value = 6742
$ws.Range("R13").Value = 'Provincia de Linares'
$ws.Range("S13").Value = 3371

# Row 14 (was row 3 data in source order)
$ws.Range("D14").Value = 44589
$ws.Range("M14").Value = 60

# Row 15 (was row 14 data in source order)
$ws.Range("D15").Value = 44211
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6000
$ws.Range("P15").Value = 6000
$ws.Range("R15").Value = 'Provincia de Curicó'
$ws.Range("S15").Value = 3000

# Row 16 (was row 18 data in source order)
$ws.Range("D16").Value = 44960
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 7000
$ws.Range("P16").Value = 7000
$ws.Range("S16").Value = 3500

# Row 17 (was row 6 data in source order)
$ws.Range("D17").Value = 44974
$ws.Range("M17").Value = 130
$ws.Range("N17").Value = 7000
$ws.Range("O17").Value = 7500
$ws.Range("P17").Value = 7269
$ws.Range("R17").Value = 'Provincia de Curicó'
$ws.Range("S17").Value = 3634

# Row 18 (was row 7 data in source order)
$ws.Range("D18").Value = 44592
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = 8000
$ws.Range("O18").Value = 8000
$ws.Range("P18").Value = 8000
$ws.Range("R18").Value = 'Provincia de Linares'
$ws.Range("S18").Value = 4000

# Row 19 (was row 13 data in source order)
$ws.Range("D19").Value = 44582
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 6000
$ws.Range("O19").Value = 6500
$ws.Range("P19").Value = 6233
$ws.Range("S19").Value = 3116
